$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AK_1.png", 0.462, 0.002, 0.605, 0.3, "Tidak Diketahui", "Salah"),
    @("AK_2.png", 0.879, 0.003, 0.823, 0.533, "Akhlak Kamiswara", "Benar"),
    @("AK_3.png", 0.833, 0.003, 0.825, 0.5669999999999999, "Akhlak Kamiswara", "Benar"),
    @("AK_4.png", 0.44, 0.001, 0.593, 0.3, "Tidak Diketahui", "Salah"),
    @("AK_5.png", 0.454, 0.001, 0.83, 0.333, "Tidak Diketahui", "Salah"),
    @("MIB_1.png", 1.623, 0.005, 0.628, 0.7, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_2.png", 1.173, 0.004, 0.725, 0.633, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_3.png", 1.744, 0.006, 0.915, 0.9330000000000001, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_4.png", 1.164, 0.004, 0.719, 0.7, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_5.png", 1.574, 0.005, 0.698, 0.7, "Muhammad Iqbal Baqi", "Benar"),
    @("AAH_1.png", 0.858, 0.003, 0.749, 0.867, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_2.png", 1.142, 0.004, 0.909, 0.967, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_3.png", 0.884, 0.003, 0.758, 0.9330000000000001, "Andrea Ayunove Hutami", "Benar"),
    @("TI_1.png", 0.961, 0.003, 0.792, 0.3, "Tidak Diketahui", "Salah"),
    @("TI_2.png", 0.945, 0.003, 0.805, 0.5, "Toni Ismail", "Benar"),
    @("TI_3.png", 0.72, 0.002, 0.901, 0.333, "Tidak Diketahui", "Salah"),
    @("TI_4.png", 0.695, 0.002, 0.764, 0.533, "Toni Ismail", "Benar"),
    @("TI_5.png", 1.032, 0.003, 0.848, 0.667, "Toni Ismail", "Benar"),
    @("RAS_1.png", 0.645, 0.002, 0.794, 0.267, "Tidak Diketahui", "Salah"),
    @("RAS_2.png", 1.026, 0.003, 0.74, 0.467, "Tidak Diketahui", "Salah"),
    @("RAS_3.png", 0.518, 0.002, 0.665, 0.367, "Tidak Diketahui", "Salah"),
    @("RAS_4.png", 1.276, 0.004, 0.618, 0.467, "Tidak Diketahui", "Salah"),
    @("RAS_5.png", 1.171, 0.004, 0.771, 0.533, "Ridha Ayu Salsabila", "Benar"),
    @("RR_1.png", 1.478, 0.005, 0.843, 0.533, "Rafiqo Rapitasari", "Benar"),
    @("RR_2.png", 1.362, 0.005, 0.858, 0.5669999999999999, "Rafiqo Rapitasari", "Benar"),
    @("RR_3.png", 1.107, 0.004, 0.503, 0.533, "Rafiqo Rapitasari", "Benar"),
    @("RR_4.png", 1.366, 0.005, 0.866, 0.5, "Rafiqo Rapitasari", "Benar"),
    @("RR_5.png", 1.479, 0.005, 0.863, 0.533, "Rafiqo Rapitasari", "Benar"),
    @("AR_1.png", 0.86, 0.003, 0.746, 0.667, "Arizli Romadhon", "Benar"),
    @("GA_1.png", 1.658, 0.005, 0.838, 0.9, "Gege Ardiyansyah", "Benar"),
    @("GA_2.png", 0.652, 0.002, 0.773, 0.233, "Tidak Diketahui", "Salah"),
    @("GA_3.png", 0.75, 0.002, 0.728, 0.3, "Tidak Diketahui", "Salah"),
    @("FY_1.png", 1.119, 0.004, 0.762, 0.233, "Tidak Diketahui", "Salah"),
    @("FY_2.png", 1.741, 0.006, 0.763, 0.333, "Tidak Diketahui", "Salah"),
    @("FY_3.png", 1.444, 0.005, 0.873, 0.267, "Tidak Diketahui", "Salah"),
    @("FY_4.png", 1.277, 0.004, 0.72, 0.267, "Tidak Diketahui", "Salah"),
    @("TO_1.png", 0.781, 0.003, 0.772, 0.533, "Tiara Oktavian", "Benar"),
    @("TO_2.png", 1.069, 0.004, 0.8129999999999999, 0.533, "Tiara Oktavian", "Benar"),
    @("TO_3.png", 0.819, 0.003, 0.8100000000000001, 0.533, "Tiara Oktavian", "Benar"),
    @("TO_4.png", 7.227, 0.024, 0.47, 0.333, "Tidak Diketahui", "Salah"),
    @("TO_5.png", 7.116, 0.023, 0.481, 0.433, "Tidak Diketahui", "Salah"),
    @("TD_1.png", 3.736, 0.012, 0.448, 0.2, "Tidak Diketahui", "Benar"),
    @("TD_2.png", 5.014, 0.017, 0.456, 0.267, "Tidak Diketahui", "Benar"),
    @("TD_3.png", 1.043, 0.003, 0.695, 0.233, "Tidak Diketahui", "Benar"),
    @("TD_4.png", 0.999, 0.003, 0.666, 0.233, "Tidak Diketahui", "Benar")
)

$nRows = $data.Count
$nCols = 7
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $data[$i][$j]
    }
}
